$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# Add the 2023 (column S) figures to the Guria region transport table.
# ------------------------------------------------------------------

# Header: 2023 column, same look as the other year headers
$ws.Range("R3").Copy($ws.Range("S3"))
$ws.Range("S3").Value = 2023

# Data rows: copy the number format from the matching source column first,
# then write the new figure (mirrors "copy a cell down, then type new value").
$ws.Range("R4").Copy($ws.Range("S4"))
$ws.Range("S4").Value = 66.400000000000006

$ws.Range("R5").Copy($ws.Range("S5"))
$ws.Range("S5").Value = 60.2

$ws.Range("R6").Copy($ws.Range("S6"))
$ws.Range("S6").Value = 1204

$ws.Range("R7").Copy($ws.Range("S7"))
$ws.Range("S7").Value = 702

$ws.Range("R8").Copy($ws.Range("S8"))
$ws.Range("S8").Value = 984.4

$ws.Range("R9").Copy($ws.Range("S9"))
$ws.Range("S9").Value = 39.4

$ws.Range("P10").Copy($ws.Range("S10"))
$ws.Range("S10").Value = 8.3000000000000007

$ws.Range("R11").Copy($ws.Range("S11"))
$ws.Range("S11").Value = 20.8

$ws.Range("R12").Copy($ws.Range("S12"))
$ws.Range("S12").Value = 1.2

$ws.Range("P13").Copy($ws.Range("S13"))
$ws.Range("S13").Value = 45.6

$ws.Range("R14").Copy($ws.Range("S14"))
$ws.Range("S14").Value = 6.2

# ------------------------------------------------------------------
# Housekeeping: a handful of empty footnote-area cells carried a
# redundant "apply number format" flavour of the same style; collapse
# them onto the plain equivalent already used by their neighbours.
# ------------------------------------------------------------------
$plain = "B18"
$dupTargets = @("C18","C19","B20","B21","E21","C22","B23","C24","E24","E26")
foreach ($t in $dupTargets) {
    $ws.Range($plain).Copy($ws.Range($t))
}

# ------------------------------------------------------------------
# Leave the freshly-entered block selected, as it would be right after
# pasting the new column in.
# ------------------------------------------------------------------
$ws.Range("S4:S14").Select()
